$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M41").Value = -3332963.8
$ws.Range("N41").Value = -7650
$ws.Range("H41").Value = 838428.4399999999
$ws.Range("L41").Value = 6770
$ws.Range("K41").Value = 3333403.8
$ws.Range("I41").Value = 3333403.8
$ws.Range("J41").Value = 6770
$ws.Range("M62").Value = -40687.875
$ws.Range("H62").Value = 47535.42
$ws.Range("I62").Value = 41311.875
$ws.Range("K62").Value = 41311.875
$ws.Range("H65").Value = 47535.42
$ws.Range("M65").Value = -203439.375
$ws.Range("I65").Value = 41311.875
$ws.Range("K65").Value = 206559.375
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80
$ws.Range("I107").Value = 2000
$ws.Range("H107").Value = 2125
$ws.Range("H137").Value = 52636268
$ws.Range("K137").Value = 333341130
$ws.Range("I137").Value = 111113710
$ws.Range("M137").Value = -333338580
$ws.Range("L138").Value = 12993.3531
$ws.Range("H138").Value = 3427.0789
$ws.Range("N138").Value = -23273.3531
$ws.Range("J138").Value = 4331.1177
$ws.Range("K141").Value = 8225.625
$ws.Range("I141").Value = 2741.875
$ws.Range("M141").Value = -3045.625
$ws.Range("H141").Value = 2698.2354

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 813.4666999999999
$ws.Range("I2").Value = 813.4666999999999
$ws.Range("M2").Value = -700.4666999999999
$ws.Range("K2").Value = 813.4666999999999
$ws.Range("H32").Value = 184291.2
$ws.Range("M32").Value = -202039.95
$ws.Range("K32").Value = 202326.95
$ws.Range("I32").Value = 202326.95
$ws.Range("N45").Value = -4476.75
$ws.Range("L45").Value = 3722.75
$ws.Range("J45").Value = 3722.75
$ws.Range("H45").Value = 3610.375
$ws.Range("H61").Value = 1753.7778
$ws.Range("J61").Value = 2733
$ws.Range("N61").Value = -3157
$ws.Range("L61").Value = 2733
$ws.Range("M61").Value = -1052.1666
$ws.Range("I61").Value = 1264.1666
$ws.Range("K61").Value = 1264.1666
$ws.Range("K116").Value = 813.4666999999999
$ws.Range("H116").Value = 813.4666999999999
$ws.Range("I116").Value = 813.4666999999999
$ws.Range("M116").Value = 1480.5333
$ws.Range("K132").Value = 8100
$ws.Range("I132").Value = 2700
$ws.Range("M132").Value = -5570
$ws.Range("H132").Value = 2700
$ws.Range("K136").Value = 3792.4998
$ws.Range("H136").Value = 1753.7778
$ws.Range("I136").Value = 1264.1666
$ws.Range("N136").Value = -13299
$ws.Range("J136").Value = 2733
$ws.Range("M136").Value = -1242.4998
$ws.Range("L136").Value = 8199

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").Value = -699.4666999999999
$ws.Range("I3").Value = 813.4666999999999
$ws.Range("H3").Value = 813.4666999999999
$ws.Range("K3").Value = 813.4666999999999
$ws.Range("J50").Value = 46999.5
$ws.Range("L50").Value = 46999.5
$ws.Range("H50").Value = 46999.5
$ws.Range("N50").Value = -48147.5
$ws.Range("I86").Value = 35715876
$ws.Range("K86").Value = 35715876
$ws.Range("M86").Value = -35714753
$ws.Range("H86").Value = 23811536
$ws.Range("H89").Value = 23811536
$ws.Range("K89").Value = 178579380
$ws.Range("I89").Value = 35715876
$ws.Range("M89").Value = -178573764
$ws.Range("M105").Value = 36.04160000000002
$ws.Range("J105").Value = 1567.9375
$ws.Range("K105").Value = 1710.9584
$ws.Range("N105").Value = -5061.9375
$ws.Range("L105").Value = 1567.9375
$ws.Range("H105").Value = 1653.75
$ws.Range("I105").Value = 1710.9584
$ws.Range("N134").Value = -21944.25
$ws.Range("J134").Value = 5624.75
$ws.Range("H134").Value = 2261.3928
$ws.Range("L134").Value = 16874.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K50").Value = 24000
$ws.Range("I50").Value = 24000
$ws.Range("H50").Value = 41000
$ws.Range("M50").Value = -23375
$ws.Range("L58").Value = 2155.25
$ws.Range("M58").Value = -1113.3636
$ws.Range("J58").Value = 2155.25
$ws.Range("N58").Value = -2561.25
$ws.Range("I58").Value = 1316.3636
$ws.Range("H58").Value = 1857.5807
$ws.Range("K58").Value = 1316.3636
$ws.Range("I86").Value = 212622
$ws.Range("K86").Value = 212622
$ws.Range("M86").Value = -211499
$ws.Range("H86").Value = 144247.17
$ws.Range("H89").Value = 144247.17
$ws.Range("K89").Value = 1063110
$ws.Range("I89").Value = 212622
$ws.Range("M89").Value = -1057494
$ws.Range("M105").Value = 727.6667
$ws.Range("J105").Value = 923
$ws.Range("K105").Value = 1019.3333
$ws.Range("N105").Value = -4417
$ws.Range("L105").Value = 923
$ws.Range("H105").Value = 1000.06665
$ws.Range("I105").Value = 1019.3333
$ws.Range("J132").Value = 4802.6665
$ws.Range("N132").Value = -19467.9995
$ws.Range("I132").Value = 3457.625
$ws.Range("L132").Value = 14407.9995
$ws.Range("K132").Value = 10372.875
$ws.Range("M132").Value = -7842.875
$ws.Range("H132").Value = 4169.706
$ws.Range("I134").Value = 2837.5833
$ws.Range("H134").Value = 3050.2334
$ws.Range("M134").Value = -5977.749899999999
$ws.Range("K134").Value = 8512.749899999999
$ws.Range("K136").Value = 3949.0908
$ws.Range("H136").Value = 1857.5807
$ws.Range("I136").Value = 1316.3636
$ws.Range("N136").Value = -11565.75
$ws.Range("J136").Value = 2155.25
$ws.Range("M136").Value = -1399.0908
$ws.Range("L136").Value = 6465.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 118.84615
$ws.Range("J38").Value = 31.75
$ws.Range("N38").Value = -789.25
$ws.Range("K38").Value = 472.66668
$ws.Range("M38").Value = -125.66668
$ws.Range("L38").Value = 95.25
$ws.Range("I38").Value = 157.55556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K2").Value = 1302.5
$ws.Range("M2").Value = -1189.5
$ws.Range("I2").Value = 1302.5
$ws.Range("N2").Value = -263.333332
$ws.Range("J2").Value = 37.333332
$ws.Range("L2").Value = 37.333332
$ws.Range("H2").Value = 957.4545000000001
$ws.Range("J132").Value = 6581.1665
$ws.Range("N132").Value = -24803.4995
$ws.Range("I132").Value = 2007600
$ws.Range("L132").Value = 19743.4995
$ws.Range("K132").Value = 6022800
$ws.Range("M132").Value = -6020270
$ws.Range("H132").Value = 916135.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I34").Value = 17499
$ws.Range("N34").Value = -47342
$ws.Range("J34").Value = 46998
$ws.Range("L34").Value = 46998
$ws.Range("H34").Value = 32248.5
$ws.Range("K34").Value = 17499
$ws.Range("M34").Value = -17327
$ws.Range("H46").Value = 2357.8518
$ws.Range("J46").Value = 2513.15
$ws.Range("L46").Value = 2513.15
$ws.Range("N46").Value = -2889.15
$ws.Range("H61").Value = 364937.25
$ws.Range("M61").Value = -377881.06
$ws.Range("I61").Value = 378083.06
$ws.Range("K61").Value = 378083.06
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("K113").Value = 378083.06
$ws.Range("H113").Value = 364937.25
$ws.Range("I113").Value = 378083.06
$ws.Range("M113").Value = -375913.06
$ws.Range("M122").Value = -5817.000100000001
$ws.Range("H122").Value = 3499.923
$ws.Range("K122").Value = 8267.000100000001
$ws.Range("I122").Value = 2755.6667
$ws.Range("N127").Value = -86236.336
$ws.Range("J127").Value = 76316.336
$ws.Range("H127").Value = 76316.336
$ws.Range("L127").Value = 76316.336
$ws.Range("H135").Value = 58666.332
$ws.Range("L135").Value = 58666.332
$ws.Range("J135").Value = 58666.332
$ws.Range("N135").Value = -68806.33199999999
$ws.Range("K136").Value = 11238
$ws.Range("H136").Value = 6758.4116
$ws.Range("I136").Value = 3746
$ws.Range("N136").Value = -33408.333
$ws.Range("J136").Value = 9436.111000000001
$ws.Range("M136").Value = -8688
$ws.Range("L136").Value = 28308.333
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K107").Value = 1715.0625
$ws.Range("M107").Value = 204.9375
$ws.Range("I107").Value = 571.6875
$ws.Range("H107").Value = 573.35297
$ws.Range("J132").Value = 4097.75
$ws.Range("N132").Value = -17353.25
$ws.Range("I132").Value = 372755.06
$ws.Range("L132").Value = 12293.25
$ws.Range("K132").Value = 1118265.18
$ws.Range("M132").Value = -1115735.18
$ws.Range("H132").Value = 288490.53

Write-Output "Applied all cell updates"
